# Change report title/header block:
#  - Row 1: title -> "Отчёт о движении продуктов" (kept merged A1:G1)
#  - Insert a new row 2 (subtitle "Период: ...", merged A2:G2, left-aligned, bordered)
#  - Insert a new blank spacer row 3 (bordered, no text)
#  - Existing header + data rows shift down by two (now rows 4-8)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert two new rows above the old header row (row 2).
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# New title text (row 1 keeps its original formatting/merge).
$ws.Range("A1:G1").Value = "Отчёт о движении продуктов"

# New subtitle row (row 2): merged, bordered, left-aligned.
$ws.Range("A2:G2").Merge()
$ws.Range("A2:G2").Value = "Период: 2023-10-14 - 2023-10-29"
$ws.Range("A2:G2").Borders.Item(1).LineStyle = 1
$ws.Range("A2:G2").Borders.Item(2).LineStyle = 1
$ws.Range("A2:G2").Borders.Item(3).LineStyle = 1
$ws.Range("A2:G2").Borders.Item(4).LineStyle = 1
$ws.Range("A2:G2").HorizontalAlignment = -4131

# New blank spacer row (row 3): bordered only, no alignment override, no content.
$ws.Range("A3:G3").Borders.Item(1).LineStyle = 1
$ws.Range("A3:G3").Borders.Item(2).LineStyle = 1
$ws.Range("A3:G3").Borders.Item(3).LineStyle = 1
$ws.Range("A3:G3").Borders.Item(4).LineStyle = 1

Write-Output "done"
